# Weekly fruit/vegetable price update: insert a new daily record as row 25
# (Vega Monumental Concepción - Poroto verde), shifting the existing rows
# 25-47 down to 26-48.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 25; this pushes the old rows 25..47 down to
# 26..48 and extends the used range to A1:R48.
$ws.Rows.Item(25).Insert()

# Populate the newly inserted row 25 with the new record's data.
$ws.Cells.Item(25, 1).Value = 11
$ws.Cells.Item(25, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(25, 3).Value = "Bíobío"
$ws.Cells.Item(25, 4).Value = 44658
$ws.Cells.Item(25, 5).Value = 8
$ws.Cells.Item(25, 6).Value = 100112031
$ws.Cells.Item(25, 7).Value = "Poroto verde"
$ws.Cells.Item(25, 8).Value = "Sin especificar"
$ws.Cells.Item(25, 9).Value = "Primera"
$ws.Cells.Item(25, 10).Value = 80
$ws.Cells.Item(25, 11).Value = 25000
$ws.Cells.Item(25, 12).Value = 25000
$ws.Cells.Item(25, 13).Value = 25000
$ws.Cells.Item(25, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(25, 15).Value = "Región Metropolitana"
$ws.Cells.Item(25, 16).Value = 1000
$ws.Cells.Item(25, 17).Value = 25
$ws.Cells.Item(25, 18).Value = "Hortaliza"
